$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 (entire rows): "Hey Java" line and the "print 'hello';" line.
# This leaves "print 'Hey';" (originally row 5) shifted up into row 3.
$ws.Rows.Item(3).Resize(2).Delete() | Out-Null

# Update the selection on the sheet to B14
$ws.Range("B14").Select() | Out-Null

# Reflect the updated workbook window geometry (best-effort; window frame
# metrics are host/session state and may not round-trip through automation).
$excel.ActiveWindow.Left = 28680
$excel.ActiveWindow.Top = -120
$excel.ActiveWindow.Width = 29040
$excel.ActiveWindow.Height = 15720
